$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "27.047.70"
Set-TextValue $ws.Range("E2") "  +0.44%  "
Set-TextValue $ws.Range("D3") "1.827.13"
Set-TextValue $ws.Range("E3") "  +0.66%  "
Set-TextValue $ws.Range("D4") "1.006"
Set-TextValue $ws.Range("E4") "  +0.74%  "
Set-TextValue $ws.Range("D5") "311.49"
Set-TextValue $ws.Range("E5") "  +0.48%  "
Set-TextValue $ws.Range("E6") "  +0.53%  "
Set-TextValue $ws.Range("D7") "0.4696"
Set-TextValue $ws.Range("E7") "  -0.64%  "
Set-TextValue $ws.Range("D8") "0.3679"
Set-TextValue $ws.Range("E8") "  -0.70%  "
Set-TextValue $ws.Range("D9") "0.07373"
Set-TextValue $ws.Range("E9") "  -0.11%  "
Set-TextValue $ws.Range("E10") "  +0.73%  "
Set-TextValue $ws.Range("D11") "20.32"
Set-TextValue $ws.Range("E11") "  -0.50%  "
Set-TextValue $ws.Range("D12") "1.833.91"
Set-TextValue $ws.Range("E12") "  -1.83%  "
Set-TextValue $ws.Range("D13") "0.07288"
Set-TextValue $ws.Range("E13") "  +3.16%  "
Set-TextValue $ws.Range("D14") "5.453"
Set-TextValue $ws.Range("E14") "  +1.87%  "
Set-TextValue $ws.Range("B15") "Litecoin"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D15") "92.24"
Set-TextValue $ws.Range("E15") "  +0.26%  "
Set-TextValue $ws.Range("B16") "Chainlink"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D16") "6.530"
Set-TextValue $ws.Range("E16") "  +0.30%  "
Set-TextValue $ws.Range("E17") "  +0.73%  "
Set-TextValue $ws.Range("D18") "0.000008762"
Set-TextValue $ws.Range("E18") "  +0.47%  "
Set-TextValue $ws.Range("E19") "  +0.43%  "
Set-TextValue $ws.Range("D20") "14.73"
Set-TextValue $ws.Range("E20") "  +0.17%  "
Set-TextValue $ws.Range("D21") "27.074.11"
Set-TextValue $ws.Range("E21") "  +0.50%  "
Set-TextValue $ws.Range("D22") "5.298"
Set-TextValue $ws.Range("E22") "  -0.74%  "
Set-TextValue $ws.Range("D23") "10.64"
Set-TextValue $ws.Range("E23") "  +0.80%  "
Set-TextValue $ws.Range("D24") "2.051.33"
Set-TextValue $ws.Range("E24") "  -2.35%  "
Set-TextValue $ws.Range("D25") "1.895"
Set-TextValue $ws.Range("E25") "  -0.03%  "
Set-TextValue $ws.Range("D26") "151.89"
Set-TextValue $ws.Range("E26") "  +0.09%  "
Set-TextValue $ws.Range("D27") "18.42"
Set-TextValue $ws.Range("E27") "  -0.07%  "
Set-TextValue $ws.Range("D28") "2.151"
Set-TextValue $ws.Range("E28") "  +1.03%  "
Set-TextValue $ws.Range("D29") "5.239"
Set-TextValue $ws.Range("E29") "  -0.97%  "
Set-TextValue $ws.Range("D30") "117.05"
Set-TextValue $ws.Range("E30") "  +1.36%  "
Set-TextValue $ws.Range("E31") "  -0.69%  "
Set-TextValue $ws.Range("D32") "0.7558"
Set-TextValue $ws.Range("E32") "  -0.20%  "
Set-TextValue $ws.Range("D33") "1.164"
Set-TextValue $ws.Range("E33") "  +0.67%  "
Set-TextValue $ws.Range("D34") "4.520"
Set-TextValue $ws.Range("E34") "  +1.09%  "
Set-TextValue $ws.Range("D35") "2.925"
Set-TextValue $ws.Range("E35") "  +0.12%  "
Set-TextValue $ws.Range("E36") "  +0.69%  "
Set-TextValue $ws.Range("E37") "  +0.58%  "
Set-TextValue $ws.Range("D38") "0.05321"
Set-TextValue $ws.Range("E38") "  +1.13%  "
Set-TextValue $ws.Range("D39") "0.01952"
Set-TextValue $ws.Range("E39") "  -0.12%  "
Set-TextValue $ws.Range("D40") "2.986"
Set-TextValue $ws.Range("E40") "  +2.34%  "
Set-TextValue $ws.Range("B41") "RenderToken"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D41") "2.399"
Set-TextValue $ws.Range("E41") "  +1.46%  "
Set-TextValue $ws.Range("B42") "FraxShare"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D42") "7.254"
Set-TextValue $ws.Range("E42") "  +0.86%  "
Set-TextValue $ws.Range("E43") "  -0.81%  "
Set-TextValue $ws.Range("D44") "0.1659"
Set-TextValue $ws.Range("E44") "  -0.14%  "
Set-TextValue $ws.Range("D45") "8.523"
Set-TextValue $ws.Range("E45") "  +0.89%  "
Set-TextValue $ws.Range("D46") "0.4928"
Set-TextValue $ws.Range("E46") "  -0.44%  "
Set-TextValue $ws.Range("D47") "10.49"
Set-TextValue $ws.Range("E47") "  +1.92%  "
Set-TextValue $ws.Range("E48") "  +0.38%  "
Set-TextValue $ws.Range("D49") "1.665"
Set-TextValue $ws.Range("E49") "  -0.48%  "
Set-TextValue $ws.Range("D50") "103.47"
Set-TextValue $ws.Range("E50") "  +0.27%  "
Set-TextValue $ws.Range("D51") "0.06306"
Set-TextValue $ws.Range("E51") "  +0.30%  "
